$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 400
$ws.Range("I16").Value = 400
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 400
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -170

$ws.Range("H33").Value = 199.71428
$ws.Range("I33").Value = 199.71428
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 199.71428
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = 29.28572

$ws.Range("H38").Value = 2731.7856
$ws.Range("I38").Value = 2537.25
$ws.Range("J38").Value = 3899
$ws.Range("K38").Value = 7611.75
$ws.Range("L38").Value = 11697
$ws.Range("M38").Value = -7239.75
$ws.Range("N38").Value = -12441

$ws.Range("H106").Value = 3821.1428
$ws.Range("I106").Value = 3821.1428
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 3821.1428
$ws.Range("L106").Value = 0
$ws.Range("M106").Value = -3190.1428

$ws.Range("H135").Value = 514205.9
$ws.Range("I135").Value = 607301.25
$ws.Range("J135").Value = 2181.5
$ws.Range("K135").Value = 5465711.25
$ws.Range("L135").Value = 19633.5
$ws.Range("M135").Value = -5463176.25
$ws.Range("N135").Value = -24703.5

$ws.Range("H138").Value = 4204.367
$ws.Range("I138").Value = 1825.6471
$ws.Range("J138").Value = 5468.0625
$ws.Range("K138").Value = 5476.9413
$ws.Range("L138").Value = 16404.1875
$ws.Range("M138").Value = -336.9412999999995
$ws.Range("N138").Value = -26684.1875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 48954.953
$ws.Range("I2").Value = 53812.684
$ws.Range("J2").Value = 2806.5
$ws.Range("K2").Value = 53812.684
$ws.Range("L2").Value = 2806.5
$ws.Range("M2").Value = -53699.684
$ws.Range("N2").Value = -3032.5

$ws.Range("H45").Value = 2660.3635
$ws.Range("I45").Value = 2114.8333
$ws.Range("J45").Value = 3315
$ws.Range("K45").Value = 2114.8333
$ws.Range("L45").Value = 3315
$ws.Range("M45").Value = -1737.8333
$ws.Range("N45").Value = -4069

$ws.Range("H102").Value = 1559.2858
$ws.Range("I102").Value = 1734.1666
$ws.Range("J102").Value = 510
$ws.Range("K102").Value = 1734.1666
$ws.Range("L102").Value = 510
$ws.Range("M102").Value = -112.1666
$ws.Range("N102").Value = -3754

$ws.Range("H116").Value = 48954.953
$ws.Range("I116").Value = 53812.684
$ws.Range("J116").Value = 2806.5
$ws.Range("K116").Value = 53812.684
$ws.Range("L116").Value = 2806.5
$ws.Range("M116").Value = -51518.684
$ws.Range("N116").Value = -7394.5

$ws.Range("H132").Value = 3024.3489
$ws.Range("I132").Value = 3082.3513
$ws.Range("J132").Value = 2666.6667
$ws.Range("K132").Value = 9247.053899999999
$ws.Range("L132").Value = 8000.000100000001
$ws.Range("M132").Value = -6717.053899999999
$ws.Range("N132").Value = -13060.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 48954.953
$ws.Range("I3").Value = 53812.684
$ws.Range("J3").Value = 2806.5
$ws.Range("K3").Value = 53812.684
$ws.Range("L3").Value = 2806.5
$ws.Range("M3").Value = -53698.684
$ws.Range("N3").Value = -3034.5

$ws.Range("H60").Value = 43949.8
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 43949.8
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 43949.8
$ws.Range("N60").Value = -45147.8

$ws.Range("H99").Value = 3298.4285
$ws.Range("I99").Value = 3298.4285
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 3298.4285
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -1800.4285

$ws.Range("H140").Value = 50000
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 50000
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 50000
$ws.Range("N140").Value = -60360

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 59606.703
$ws.Range("I31").Value = 2749.5
$ws.Range("J31").Value = 64155.28
$ws.Range("K31").Value = 2749.5
$ws.Range("L31").Value = 64155.28
$ws.Range("M31").Value = -2454.5
$ws.Range("N31").Value = -64745.28

$ws.Range("H34").Value = 59606.703
$ws.Range("I34").Value = 2749.5
$ws.Range("J34").Value = 64155.28
$ws.Range("K34").Value = 2749.5
$ws.Range("L34").Value = 64155.28
$ws.Range("M34").Value = -2547.5
$ws.Range("N34").Value = -64559.28

$ws.Range("H62").Value = 5280.2
$ws.Range("I62").Value = 3500.3333
$ws.Range("J62").Value = 7950
$ws.Range("K62").Value = 3500.3333
$ws.Range("L62").Value = 7950
$ws.Range("M62").Value = -2876.3333
$ws.Range("N62").Value = -9198

$ws.Range("H65").Value = 5280.2
$ws.Range("I65").Value = 3500.3333
$ws.Range("J65").Value = 7950
$ws.Range("K65").Value = 17501.6665
$ws.Range("L65").Value = 39750
$ws.Range("M65").Value = -14381.6665
$ws.Range("N65").Value = -45990

$ws.Range("H68").Value = 107999.2
$ws.Range("I68").Value = 99998.5
$ws.Range("J68").Value = 113333
$ws.Range("K68").Value = 99998.5
$ws.Range("L68").Value = 113333
$ws.Range("M68").Value = -99249.5
$ws.Range("N68").Value = -114831

$ws.Range("H71").Value = 107999.2
$ws.Range("I71").Value = 99998.5
$ws.Range("J71").Value = 113333
$ws.Range("K71").Value = 299995.5
$ws.Range("L71").Value = 339999
$ws.Range("M71").Value = -296251.5
$ws.Range("N71").Value = -347487

$ws.Range("H132").Value = 2900
$ws.Range("I132").Value = 2900
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 8700
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -6170

$ws.Range("H134").Value = 229406.48
$ws.Range("I134").Value = 2243.5366
$ws.Range("J134").Value = 3333966.8
$ws.Range("K134").Value = 6730.6098
$ws.Range("L134").Value = 10001900.4
$ws.Range("M134").Value = -4195.6098
$ws.Range("N134").Value = -10006970.4

$ws.Range("H139").Value = 99690
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 99690
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 99690
$ws.Range("N139").Value = -109970

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 962645.3
$ws.Range("I4").Value = 1053798.2
$ws.Range("J4").Value = 251652
$ws.Range("K4").Value = 3161394.6
$ws.Range("L4").Value = 754956
$ws.Range("M4").Value = -3161282.6
$ws.Range("N4").Value = -755180

$ws.Range("H51").Value = 1366.3334
$ws.Range("I51").Value = 200
$ws.Range("J51").Value = 1949.5
$ws.Range("K51").Value = 600
$ws.Range("L51").Value = 5848.5
$ws.Range("M51").Value = -140
$ws.Range("N51").Value = -6768.5

$ws.Range("H107").Value = 103417.65
$ws.Range("I107").Value = 1023.4545
$ws.Range("J107").Value = 228566.11
$ws.Range("K107").Value = 3070.3635
$ws.Range("L107").Value = 685698.33
$ws.Range("M107").Value = -1150.3635
$ws.Range("N107").Value = -689538.33

$ws.Range("H113").Value = 1765254.4
$ws.Range("I113").Value = 9260484
$ws.Range("J113").Value = 1671
$ws.Range("K113").Value = 27781452
$ws.Range("L113").Value = 5013
$ws.Range("M113").Value = -27779282
$ws.Range("N113").Value = -9353

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 3227.318
$ws.Range("I113").Value = 2068
$ws.Range("J113").Value = 4386.636
$ws.Range("K113").Value = 2068
$ws.Range("L113").Value = 4386.636
$ws.Range("M113").Value = 102
$ws.Range("N113").Value = -8726.636

$ws.Range("H132").Value = 72666.664
$ws.Range("I132").Value = 9187.75
$ws.Range("J132").Value = 145214
$ws.Range("K132").Value = 27563.25
$ws.Range("L132").Value = 435642
$ws.Range("M132").Value = -25033.25
$ws.Range("N132").Value = -440702

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 5000
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 5000
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 5000
$ws.Range("N4").Value = -5226

$ws.Range("H16").Value = 166667940
$ws.Range("I16").Value = 333334200
$ws.Range("J16").Value = 1666.3334
$ws.Range("K16").Value = 333334200
$ws.Range("L16").Value = 1666.3334
$ws.Range("M16").Value = -333334030
$ws.Range("N16").Value = -2006.3334

$ws.Range("H28").Value = 5000
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 5000
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 5000
$ws.Range("N28").Value = -5464

$ws.Range("H37").Value = 5000
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 5000
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 5000
$ws.Range("N37").Value = -5214

$ws.Range("H68").Value = 252063.25
$ws.Range("I68").Value = 1200
$ws.Range("J68").Value = 335684.34
$ws.Range("K68").Value = 1200
$ws.Range("L68").Value = 335684.34
$ws.Range("M68").Value = -451
$ws.Range("N68").Value = -337182.34

$ws.Range("H71").Value = 252063.25
$ws.Range("I71").Value = 1200
$ws.Range("J71").Value = 335684.34
$ws.Range("K71").Value = 6000
$ws.Range("L71").Value = 1678421.7
$ws.Range("M71").Value = -2256
$ws.Range("N71").Value = -1685909.7

$ws.Range("H100").Value = 3333
$ws.Range("I100").Value = 2999.5
$ws.Range("J100").Value = 4000
$ws.Range("K100").Value = 2999.5
$ws.Range("L100").Value = 4000
$ws.Range("M100").Value = -2458.5
$ws.Range("N100").Value = -5082

$ws.Range("H132").Value = 7247.2354
$ws.Range("I132").Value = 6168.125
$ws.Range("J132").Value = 8206.444
$ws.Range("K132").Value = 18504.375
$ws.Range("L132").Value = 24619.332
$ws.Range("M132").Value = -15974.375
$ws.Range("N132").Value = -29679.332

$ws.Range("H133").Value = 53789.9
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 53789.9
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 53789.9
$ws.Range("N133").Value = -58849.9

$ws.Range("H136").Value = 561477.75
$ws.Range("I136").Value = 719649.8
$ws.Range("J136").Value = 7875.5
$ws.Range("K136").Value = 2158949.4
$ws.Range("L136").Value = 23626.5
$ws.Range("M136").Value = -2156399.4
$ws.Range("N136").Value = -28726.5

$ws.Range("H137").Value = 53750
$ws.Range("I137").Value = 65000

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 82830.62
$ws.Range("I62").Value = 204359.6
$ws.Range("J62").Value = 6875
$ws.Range("K62").Value = 204359.6
$ws.Range("L62").Value = 6875
$ws.Range("M62").Value = -203735.6
$ws.Range("N62").Value = -8123

$ws.Range("H65").Value = 82830.62
$ws.Range("I65").Value = 204359.6
$ws.Range("J65").Value = 6875
$ws.Range("K65").Value = 1021798
$ws.Range("L65").Value = 34375
$ws.Range("M65").Value = -1018678
$ws.Range("N65").Value = -40615

$ws.Range("H100").Value = 943.25
$ws.Range("I100").Value = 1007.6667
$ws.Range("J100").Value = 750
$ws.Range("K100").Value = 2015.3334
$ws.Range("L100").Value = 1500
$ws.Range("M100").Value = -1474.3334
$ws.Range("N100").Value = -2582

$ws.Range("H136").Value = 9868442
$ws.Range("I136").Value = 11840895
$ws.Range("J136").Value = 334917.34
$ws.Range("K136").Value = 35522685
$ws.Range("L136").Value = 1004752.02
$ws.Range("M136").Value = -35520135
$ws.Range("N136").Value = -1009852.02
